# Adds results generated "by Google Colab pro plus":
#   - a small 0..5 header row (row 8, columns B:G)
#   - a 6x6 EXP() lookup table (rows 9:14, columns B:G) referencing the
#     original A1:F6 data block, with a 0..5 index column (A9:A14)
#   - an integer number format on the header/index cells
#   - conditional formatting (red font) on the EXP table when value >= 0.3
#   - best-fit-ish column widths for columns A:C
#   - updates the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Integer number format used for the header row and the index column ---
$intFmt = "0_);[Red](0)"

# --- Row 8: header values 0..5 in columns B..G ---
for ($c = 0; $c -lt 6; $c++) {
    $cell = $ws.Cells.Item(8, $c + 2)
    $cell.NumberFormat = $intFmt
    $cell.Value = $c
}

# --- Rows 9..14: index column (A) + EXP() table (B..G) ---
$srcCols = @("A", "B", "C", "D", "E", "F")

for ($r = 0; $r -lt 6; $r++) {
    $destRow = $r + 9
    $srcRow = $r + 1

    $idxCell = $ws.Cells.Item($destRow, 1)
    $idxCell.NumberFormat = $intFmt
    $idxCell.Value = $r

    for ($c = 0; $c -lt 6; $c++) {
        $destCol = $c + 2
        $ws.Cells.Item($destRow, $destCol).Formula = "=EXP(" + $srcCols[$c] + $srcRow + ")"
    }
}

# --- Conditional formatting: red font when value >= 0.3 ---
$rng = $ws.Range("B9:G14")
$fc = $rng.FormatConditions.Add(1, 7, "0.3")
$fc.Font.Color = 255

# --- Column widths for A:C ---
$ws.Range("A1:C1").ColumnWidth = 8.449776785714286

# --- Update selection ---
$null = $ws.Range("E18").Select()
